# edit.ps1 - applies the "ppt update, mobile responsive" change set
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 ("Concept"): expand "User story" bullet into full user story
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Para = $s3Body.Paragraphs(3,1)
$s3Para.Text = "User story: As a new developer, I want an app which can be used to help “break the ice”, so that I can get to know my fellow developers better."

# ---------------------------------------------------------------------
# Slide 4 ("Process"): add technology bullets + task breakdown bullets
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange

# --- insert the technology sub-bullets right before "Breakdown of tasks and roles"
$breakdownPara = $s4Body.Paragraphs(2,1)
$techBlock = "passport.js" + [char]13 + "express" + [char]13 + "handlebars" + [char]13 + "Sequelize " + [char]13 + [char]13
$null = $breakdownPara.InsertBefore($techBlock)

$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
for ($i = 2; $i -le 6; $i++) {
    $para = $s4Body.Paragraphs($i,1)
    $para.IndentLevel = 2
}

# split "Sequelize " into its own run plus a trailing space run
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$seqPara = $s4Body.Paragraphs(5,1)
$seqSpace = $seqPara.Characters(10,1)
$seqSpace.Text = $seqSpace.Text

# --- insert the task/role breakdown sub-bullets right after "Breakdown of tasks and roles"
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$challengesPara = $s4Body.Paragraphs(8,1)
$roleBlock = "Analee – routes, passport.js" + [char]13 + "Keita – models, PowerPoint" + [char]13 + "Kevin – routes, html, " + [char]13 + "Laura – models, CSS" + [char]13
$null = $challengesPara.InsertBefore($roleBlock)

$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
for ($i = 8; $i -le 11; $i++) {
    $para = $s4Body.Paragraphs($i,1)
    $para.IndentLevel = 2
}

# split "Analee – routes, passport.js" into 3 runs: "Analee" / " – routes, " / "passport.js"
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$analeePara = $s4Body.Paragraphs(8,1)
$midRun = $analeePara.Characters(7,11)
$midRun.Text = $midRun.Text

# ---------------------------------------------------------------------
# Slide 5 ("Demo"): fill in the empty content placeholder
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5Body = $s5.Shapes.Item(2).TextFrame.TextRange
$s5Body.Text = "Demo video here"

# ---------------------------------------------------------------------
# Slide 6 ("Directions for Future Development"): add the 4 bullet points
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6Body = $s6.Shapes.Item(2).TextFrame.TextRange
$s6Body.Text = "Have user to create questions" + [char]13 + "Expand answer variety" + [char]13 + "Create a game to guess who is the developer based on the answers provided" + [char]13 + "Track score"

# ---------------------------------------------------------------------
# Slide 7 ("Links"): mark the "GitHub repo" run dirty
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7Body = $s7.Shapes.Item(2).TextFrame.TextRange
$repoPara = $s7Body.Paragraphs(2,1)
$repoPara.Text = $repoPara.Text
